$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44523.81192129629
$ws.Range("C5").Value = 44523.81467592593
$ws.Range("D5").Value = "IP Address"
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 237
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 44523.81467592593
$ws.Range("I5").Value = "1dabec"
$ws.Range("J5").Value = "ebola %>%`n  pivot_longer (``289``:last_col(), names_to = ""day"", values_to = ""cases"") %>%`n  na.omit()"

$ws.Range("B5:C5").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("H5").NumberFormat = $ws.Range("H2").NumberFormat
$ws.Rows(5).AutoFit()

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44518.81696759259
$ws.Range("C6").Value = 44518.85575231482
$ws.Range("D6").Value = "IP Address"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 3350
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 44523.87966435185
$ws.Range("I6").Value = "2hagra"

$ws.Range("B6:C6").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("H6").NumberFormat = $ws.Range("H2").NumberFormat

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44523.1821412037
$ws.Range("C7").Value = 44523.18387731482
$ws.Range("D7").Value = "Spam"
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 149
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 44523.87975694444
$ws.Range("I7").Value = "2nesch"

$ws.Range("B7:C7").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("H7").NumberFormat = $ws.Range("H2").NumberFormat
